$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 4).Value = Get-Date -Year 2023 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100104
$ws.Cells.Item($row, 8).Value = "Frutos de pepita"
$ws.Cells.Item($row, 9).Value = 100104003
$ws.Cells.Item($row, 10).Value = "Membrillo"
$ws.Cells.Item($row, 11).Value = "Champion"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 50
$ws.Cells.Item($row, 14).Value = 12000
$ws.Cells.Item($row, 15).Value = 12000
$ws.Cells.Item($row, 16).Value = 12000
$ws.Cells.Item($row, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 667
$ws.Cells.Item($row, 20).Value = 18
